$d = $word.ActiveDocument

# Update the date heading paragraph
$d.Paragraphs.Item(1).Range.Text = "2024-09-11 Wednesday"

# Update each answer cell in the table by (row, column) position
# to avoid ambiguity from duplicate old values appearing more than once.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "25+67=92"
$t.Cell(1, 2).Range.Text = "38-9=29"
$t.Cell(1, 3).Range.Text = "5+38=43"
$t.Cell(1, 4).Range.Text = "26+45=71"
$t.Cell(1, 5).Range.Text = "85-69=16"
$t.Cell(2, 1).Range.Text = "88+3=91"
$t.Cell(2, 2).Range.Text = "24+29=53"
$t.Cell(2, 3).Range.Text = "57+29=86"
$t.Cell(2, 4).Range.Text = "70-41=29"
$t.Cell(2, 5).Range.Text = "17+78=95"
$t.Cell(3, 1).Range.Text = "22-14=8"
$t.Cell(3, 2).Range.Text = "70-52=18"
$t.Cell(3, 3).Range.Text = "51-22=29"
$t.Cell(3, 4).Range.Text = "91-3=88"
$t.Cell(3, 5).Range.Text = "22-7=15"
$t.Cell(4, 1).Range.Text = "45-29=16"
$t.Cell(4, 2).Range.Text = "93-89=4"
$t.Cell(4, 3).Range.Text = "82-54=28"
$t.Cell(4, 4).Range.Text = "35-27=8"
$t.Cell(4, 5).Range.Text = "60-48=12"
$t.Cell(5, 1).Range.Text = "18+43=61"
$t.Cell(5, 2).Range.Text = "5+38=43"
$t.Cell(5, 3).Range.Text = "29+63=92"
$t.Cell(5, 4).Range.Text = "41-8=33"
$t.Cell(5, 5).Range.Text = "40-11=29"
$t.Cell(6, 1).Range.Text = "78+14=92"
$t.Cell(6, 2).Range.Text = "52-16=36"
$t.Cell(6, 3).Range.Text = "31-29=2"
$t.Cell(6, 4).Range.Text = "44+8=52"
$t.Cell(6, 5).Range.Text = "20-3=17"
$t.Cell(7, 1).Range.Text = "9+12=21"
$t.Cell(7, 2).Range.Text = "81-64=17"
$t.Cell(7, 3).Range.Text = "60-1=59"
$t.Cell(7, 4).Range.Text = "19+46=65"
$t.Cell(7, 5).Range.Text = "15+16=31"
$t.Cell(8, 1).Range.Text = "36+37=73"
$t.Cell(8, 2).Range.Text = "39+47=86"
$t.Cell(8, 3).Range.Text = "35+16=51"
$t.Cell(8, 4).Range.Text = "26+37=63"
$t.Cell(8, 5).Range.Text = "28+7=35"
$t.Cell(9, 1).Range.Text = "87-69=18"
$t.Cell(9, 2).Range.Text = "40-28=12"
$t.Cell(9, 3).Range.Text = "80-79=1"
$t.Cell(9, 4).Range.Text = "27+6=33"
$t.Cell(9, 5).Range.Text = "58+14=72"
$t.Cell(10, 1).Range.Text = "42-27=15"
$t.Cell(10, 2).Range.Text = "80-29=51"
$t.Cell(10, 3).Range.Text = "54-36=18"
$t.Cell(10, 4).Range.Text = "54-48=6"
$t.Cell(10, 5).Range.Text = "58+5=63"
$t.Cell(11, 1).Range.Text = "29+38=67"
$t.Cell(11, 2).Range.Text = "82-43=39"
$t.Cell(11, 3).Range.Text = "92-19=73"
$t.Cell(11, 4).Range.Text = "18+64=82"
$t.Cell(11, 5).Range.Text = "27+16=43"
$t.Cell(12, 1).Range.Text = "44-18=26"
$t.Cell(12, 2).Range.Text = "90-74=16"
$t.Cell(12, 3).Range.Text = "67+25=92"
$t.Cell(12, 4).Range.Text = "72-47=25"
$t.Cell(12, 5).Range.Text = "96-57=39"
$t.Cell(13, 1).Range.Text = "8+79=87"
$t.Cell(13, 2).Range.Text = "46-29=17"
$t.Cell(13, 3).Range.Text = "45+36=81"
$t.Cell(13, 4).Range.Text = "7+39=46"
$t.Cell(13, 5).Range.Text = "15+58=73"
$t.Cell(14, 1).Range.Text = "74-7=67"
$t.Cell(14, 2).Range.Text = "8+7=15"
$t.Cell(14, 3).Range.Text = "26+68=94"
$t.Cell(14, 4).Range.Text = "34+57=91"
$t.Cell(14, 5).Range.Text = "90-3=87"
$t.Cell(15, 1).Range.Text = "39+36=75"
$t.Cell(15, 2).Range.Text = "40-22=18"
$t.Cell(15, 3).Range.Text = "92-24=68"
$t.Cell(15, 4).Range.Text = "95-78=17"
$t.Cell(15, 5).Range.Text = "94-6=88"
$t.Cell(16, 1).Range.Text = "34+7=41"
$t.Cell(16, 2).Range.Text = "57+19=76"
$t.Cell(16, 3).Range.Text = "3+19=22"
$t.Cell(16, 4).Range.Text = "63-35=28"
$t.Cell(16, 5).Range.Text = "8+39=47"
$t.Cell(17, 1).Range.Text = "63-17=46"
$t.Cell(17, 2).Range.Text = "35+46=81"
$t.Cell(17, 3).Range.Text = "35+39=74"
$t.Cell(17, 4).Range.Text = "38+46=84"
$t.Cell(17, 5).Range.Text = "55-17=38"
$t.Cell(18, 1).Range.Text = "81-57=24"
$t.Cell(18, 2).Range.Text = "86+6=92"
$t.Cell(18, 3).Range.Text = "81-62=19"
$t.Cell(18, 4).Range.Text = "9+19=28"
$t.Cell(18, 5).Range.Text = "55-37=18"
$t.Cell(19, 1).Range.Text = "61-45=16"
$t.Cell(19, 2).Range.Text = "43-6=37"
$t.Cell(19, 3).Range.Text = "46+37=83"
$t.Cell(19, 4).Range.Text = "7+26=33"
$t.Cell(19, 5).Range.Text = "8+43=51"
$t.Cell(20, 1).Range.Text = "82-7=75"
$t.Cell(20, 2).Range.Text = "92-8=84"
$t.Cell(20, 3).Range.Text = "49+14=63"
$t.Cell(20, 4).Range.Text = "93-66=27"
$t.Cell(20, 5).Range.Text = "50-29=21"
